# Auto-generated from diff: Famfrit_Profits (Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 163.38461
$ws.Range("I11").Value = 163.38461
$ws.Range("K11").Value = 163.38461
$ws.Range("M11").Value = -23.38461000000001
$ws.Range("H107").Value = 2565
$ws.Range("I107").Value = 2218.4
$ws.Range("J107").Value = 2998.25
$ws.Range("K107").Value = 2218.4
$ws.Range("L107").Value = 2998.25
$ws.Range("M107").Value = -298.4000000000001
$ws.Range("N107").Value = -6838.25
$ws.Range("H135").Value = 15625913
$ws.Range("I135").Value = 964.63635
$ws.Range("K135").Value = 8681.727150000001
$ws.Range("M135").Value = -6146.727150000001
$ws.Range("H137").Value = 4550.4614
$ws.Range("I137").Value = 1518.3478
$ws.Range("J137").Value = 27796.666
$ws.Range("K137").Value = 4555.0434
$ws.Range("L137").Value = 83389.99800000001
$ws.Range("M137").Value = -2005.0434
$ws.Range("N137").Value = -88489.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6451.4287
$ws.Range("I32").Value = 5626.091
$ws.Range("K32").Value = 5626.091
$ws.Range("M32").Value = -5339.091
$ws.Range("H37").Value = 39939
$ws.Range("J37").Value = 39939
$ws.Range("L37").Value = 39939
$ws.Range("N37").Value = -40485
$ws.Range("H45").Value = 3469
$ws.Range("I45").Value = 2855.2222
$ws.Range("J45").Value = 4850
$ws.Range("K45").Value = 2855.2222
$ws.Range("L45").Value = 4850
$ws.Range("M45").Value = -2478.2222
$ws.Range("N45").Value = -5604
$ws.Range("H61").Value = 23811618
$ws.Range("I61").Value = 29413694
$ws.Range("K61").Value = 29413694
$ws.Range("M61").Value = -29413482
$ws.Range("H122").Value = 9261126
$ws.Range("I122").Value = 926.3913
$ws.Range("K122").Value = 2779.1739
$ws.Range("M122").Value = -329.1738999999998
$ws.Range("H136").Value = 23811618
$ws.Range("I136").Value = 29413694
$ws.Range("K136").Value = 88241082
$ws.Range("M136").Value = -88238532

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3576.55
$ws.Range("I134").Value = 3585.0557
$ws.Range("K134").Value = 10755.1671
$ws.Range("M134").Value = -8220.167099999999
$ws.Range("H140").Value = 197999.2
$ws.Range("J140").Value = 197999.2
$ws.Range("L140").Value = 197999.2
$ws.Range("N140").Value = -208359.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12505371
$ws.Range("J31").Value = 35723772
$ws.Range("L31").Value = 35723772
$ws.Range("N31").Value = -35724362
$ws.Range("H34").Value = 12505371
$ws.Range("J34").Value = 35723772
$ws.Range("L34").Value = 35723772
$ws.Range("N34").Value = -35724176
$ws.Range("H58").Value = 1526.75
$ws.Range("I58").Value = 1526.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1526.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1323.75
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 14016.714
$ws.Range("J86").Value = 12777.5
$ws.Range("L86").Value = 12777.5
$ws.Range("N86").Value = -15023.5
$ws.Range("H89").Value = 14016.714
$ws.Range("J89").Value = 12777.5
$ws.Range("L89").Value = 63887.5
$ws.Range("N89").Value = -75119.5
$ws.Range("H107").Value = 993.61536
$ws.Range("I107").Value = 950.1111
$ws.Range("J107").Value = 1091.5
$ws.Range("K107").Value = 950.1111
$ws.Range("L107").Value = 1091.5
$ws.Range("M107").Value = 969.8889
$ws.Range("N107").Value = -4931.5
$ws.Range("H136").Value = 1526.75
$ws.Range("I136").Value = 1526.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4580.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2030.25
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 115136.14
$ws.Range("I141").Value = 27000
$ws.Range("J141").Value = 129825.5
$ws.Range("K141").Value = 27000
$ws.Range("L141").Value = 129825.5
$ws.Range("M141").Value = -21820
$ws.Range("N141").Value = -140185.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 196015
$ws.Range("I128").Value = 196015
$ws.Range("K128").Value = 588045
$ws.Range("M128").Value = -583065
$ws.Range("H137").Value = 2443.1428
$ws.Range("J137").Value = 2600.3333
$ws.Range("L137").Value = 7800.999899999999
$ws.Range("N137").Value = -18000.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6980.273
$ws.Range("I70").Value = 5664.6665
$ws.Range("J70").Value = 7473.625
$ws.Range("K70").Value = 5664.6665
$ws.Range("L70").Value = 7473.625
$ws.Range("M70").Value = -5394.6665
$ws.Range("N70").Value = -8013.625
$ws.Range("H73").Value = 6980.273
$ws.Range("I73").Value = 5664.6665
$ws.Range("J73").Value = 7473.625
$ws.Range("K73").Value = 5664.6665
$ws.Range("L73").Value = 7473.625
$ws.Range("M73").Value = -4728.6665
$ws.Range("N73").Value = -9345.625
$ws.Range("H102").Value = 3356.4
$ws.Range("I102").Value = 1763.3334
$ws.Range("K102").Value = 1763.3334
$ws.Range("M102").Value = -141.3334
$ws.Range("H122").Value = 16668824
$ws.Range("I122").Value = 1962.2727
$ws.Range("J122").Value = 62502692
$ws.Range("K122").Value = 5886.8181
$ws.Range("L122").Value = 187508076
$ws.Range("M122").Value = -3436.8181
$ws.Range("N122").Value = -187512976

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6824
$ws.Range("I40").Value = 5765.6665
$ws.Range("J40").Value = 9999
$ws.Range("K40").Value = 5765.6665
$ws.Range("L40").Value = 9999
$ws.Range("M40").Value = -5629.6665
$ws.Range("N40").Value = -10271
$ws.Range("H61").Value = 9642.888999999999
$ws.Range("I61").Value = 8465.333000000001
$ws.Range("K61").Value = 8465.333000000001
$ws.Range("M61").Value = -8263.333000000001
$ws.Range("H113").Value = 9642.888999999999
$ws.Range("I113").Value = 8465.333000000001
$ws.Range("K113").Value = 8465.333000000001
$ws.Range("M113").Value = -6295.333000000001
$ws.Range("H122").Value = 3575492.8
$ws.Range("I122").Value = 3775.926
$ws.Range("K122").Value = 11327.778
$ws.Range("M122").Value = -8877.778
$ws.Range("H131").Value = 63249.2
$ws.Range("J131").Value = 70316.664
$ws.Range("L131").Value = 70316.664
$ws.Range("N131").Value = -80396.664
$ws.Range("H132").Value = 3444
$ws.Range("I132").Value = 2426.0908
$ws.Range("K132").Value = 7278.2724
$ws.Range("M132").Value = -4748.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5938.4
$ws.Range("J96").Value = 3235.5
$ws.Range("L96").Value = 3235.5
$ws.Range("N96").Value = -5981.5
$ws.Range("H107").Value = 2133.25
$ws.Range("I107").Value = 1349.875
$ws.Range("J107").Value = 2655.5
$ws.Range("K107").Value = 4049.625
$ws.Range("L107").Value = 7966.5
$ws.Range("M107").Value = -2129.625
$ws.Range("N107").Value = -11806.5
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
